# Update "想去人数" (column F) counts across the three sheets that carry
# per-event data. The fourth sheet ("本地生活") has no changes.

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 1127
$ws.Range("F7").Value  = 601
$ws.Range("F8").Value  = 1538
$ws.Range("F11").Value = 1441
$ws.Range("F13").Value = 597
$ws.Range("F14").Value = 1743
$ws.Range("F15").Value = 1793
$ws.Range("F18").Value = 1458
$ws.Range("F19").Value = 285
$ws.Range("F22").Value = 1195
$ws.Range("F23").Value = 396
$ws.Range("F24").Value = 439
$ws.Range("F26").Value = 4725
$ws.Range("F29").Value = 1623
$ws.Range("F31").Value = 112

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 63

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 63
$ws.Range("F12").Value = 1127
$ws.Range("F15").Value = 601
$ws.Range("F16").Value = 1538
$ws.Range("F20").Value = 1441
$ws.Range("F22").Value = 597
$ws.Range("F23").Value = 1743
$ws.Range("F24").Value = 1793
$ws.Range("F27").Value = 1458
$ws.Range("F28").Value = 285
$ws.Range("F33").Value = 1195
$ws.Range("F34").Value = 396
$ws.Range("F35").Value = 439
$ws.Range("F37").Value = 4725
$ws.Range("F40").Value = 1623
$ws.Range("F44").Value = 112

$wb.Save()
